$wb = $excel.ActiveWorkbook

# --- term metadata bump: 1.0.0 -> 1.1.0, refreshed publish date ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.1.0"
$wsMeta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# --- fix the header/data styles so the wrap-text alignment is actually applied ---
$wsMeta.Range("A1:B14").WrapText = $true

$wsInclude = $wb.Worksheets.Item("Include from FFB")
$wsInclude.Range("A1:C2").WrapText = $true
$wsInclude.Range("A3:B4").WrapText = $true
